$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure cells keep text formatting (avoid Excel auto-converting numeric-looking strings)
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '35.340.08'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +0.03%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.913.98'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +0.36%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.12%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.724'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +8.29%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '255.08'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +3.86%  '
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -0.11%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '40.79'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -1.22%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.371'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +6.12%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '52.87'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +0.17%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0762'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +6.13%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0988'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -0.74%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '2.188.24'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +0.15%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '12.85'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +6.11%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.726'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +3.99%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '4.99'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +2.54%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.880.60'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -1.90%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '35.297.77'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -0.05%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '74.60'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +2.73%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.0₃0856'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +3.33%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '244.41'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +2.02%  '
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +4.96%  '
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +5.87%  '
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +0.01%  '
$ws.Range('B25').NumberFormat = '@'
$ws.Range('B25').Value = 'Toncoin'
$ws.Range('C25').NumberFormat = '@'
$ws.Range('C25').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.40'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +4.35%  '
$ws.Range('B26').NumberFormat = '@'
$ws.Range('B26').Value = 'PancakeSwap'
$ws.Range('C26').NumberFormat = '@'
$ws.Range('C26').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.44'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +4.91%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '167.14'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -1.76%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.68'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +2.71%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '18.81'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +2.11%  '
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +4.62%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.128.90'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +19.46%  '
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +5.46%  '
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +13.94%  '
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +23.05%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0590'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +4.44%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '4.24'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +3.37%  '
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -1.05%  '
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -2.16%  '
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +0.34%  '
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +5.27%  '
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +5.38%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '97.09'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +7.97%  '
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +1.58%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0647'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +0.50%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.337.87'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.44'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +2.14%  '
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +1.01%  '
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +3.01%  '
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -0.33%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '45.26'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -4.61%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '11.85'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +9.48%  '
